$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 491, pushing existing rows 491.. down by 2.
$ws.Rows.Item(491).Insert()
$ws.Rows.Item(491).Insert()

# New weekly entry: row 491 = "1a (cosecha)", row 492 = "2a (cosecha)"
$ws.Cells.Item(491, 1).Value = 8
$ws.Cells.Item(491, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(491, 3).Value = "Coquimbo"
$ws.Cells.Item(491, 4).Value = 44641
$ws.Cells.Item(491, 5).Value = 4
$ws.Cells.Item(491, 6).Value = 100112045
$ws.Cells.Item(491, 7).Value = "Zapallo"
$ws.Cells.Item(491, 8).Value = "Camote"
$ws.Cells.Item(491, 9).Value = "1a (cosecha)"
$ws.Cells.Item(491, 10).Value = 1480
$ws.Cells.Item(491, 11).Value = 480
$ws.Cells.Item(491, 12).Value = 500
$ws.Cells.Item(491, 13).Value = 490
$ws.Cells.Item(491, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(491, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(491, 16).Value = 490
$ws.Cells.Item(491, 17).Value = 1
$ws.Cells.Item(491, 18).Value = "Hortaliza"

$ws.Cells.Item(492, 1).Value = 8
$ws.Cells.Item(492, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(492, 3).Value = "Coquimbo"
$ws.Cells.Item(492, 4).Value = 44641
$ws.Cells.Item(492, 5).Value = 4
$ws.Cells.Item(492, 6).Value = 100112045
$ws.Cells.Item(492, 7).Value = "Zapallo"
$ws.Cells.Item(492, 8).Value = "Camote"
$ws.Cells.Item(492, 9).Value = "2a (cosecha)"
$ws.Cells.Item(492, 10).Value = 800
$ws.Cells.Item(492, 11).Value = 380
$ws.Cells.Item(492, 12).Value = 400
$ws.Cells.Item(492, 13).Value = 390
$ws.Cells.Item(492, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(492, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(492, 16).Value = 390
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"

# Keep the date columns formatted like the rest of column D.
$ws.Range("D491:D492").NumberFormat = $ws.Range("D493").NumberFormat
